# Code changes to Dashboard, My Vehicles and Home Page
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Home_Page row (row 4): Runmode changes from "N" to "Y"
$ws.Range("C4").Value = "Y"

# MyVehicles row (row 9): Done column set to "Done", Runmode changes from "Y" to "N"
$ws.Range("B9").Value = "Done"
$ws.Range("C9").Value = "N"

# Update the active selection / cursor position shown in the sheet view
$ws.Range("F11").Select() | Out-Null
